$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Value = "'62.544.51"
$c.Style = "Normal"

$c = $ws.Cells.Item(2, 5)
$c.Value = "'  -1.46%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.Value = "'3.012.11"
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 5)
$c.Value = "'  -1.77%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 5)
$c.Value = "'  +0.06%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.Value = "'585.37"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 5)
$c.Value = "'  -1.31%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 4)
$c.Value = "'146.46"
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 5)
$c.Value = "'  -5.25%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 5)
$c.Value = "'  +0.10%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 5)
$c.Value = "'  -2.47%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 4)
$c.Value = "'3.008.51"
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 5)
$c.Value = "'  -1.91%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 5)
$c.Value = "'  -4.34%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.Value = "'5.81"
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 5)
$c.Value = "'  -0.59%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.461"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 5)
$c.Value = "'  +2.19%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 5)
$c.Value = "'  -3.23%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.Value = "'34.81"
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 5)
$c.Value = "'  -5.63%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 5)
$c.Value = "'  +2.19%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.Value = "'3.507.85"
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 5)
$c.Value = "'  -1.80%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 4)
$c.Value = "'7.11"
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 5)
$c.Value = "'  -1.06%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 4)
$c.Value = "'62.515.74"
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 5)
$c.Value = "'  -1.40%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 4)
$c.Value = "'3.011.14"
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 5)
$c.Value = "'  -1.82%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 4)
$c.Value = "'459.21"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 5)
$c.Value = "'  -6.32%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.Value = "'14.00"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 5)
$c.Value = "'  -2.91%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.Value = "'0.690"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 5)
$c.Value = "'  -2.46%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.Value = "'7.44"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 5)
$c.Value = "'  -1.66%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 4)
$c.Value = "'81.80"
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.Value = "'  -0.26%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 2)
$c.Value = "'InternetComputer(DFINITY)"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 3)
$c.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 4)
$c.Value = "'12.38"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 5)
$c.Value = "'  -4.04%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 2)
$c.Value = "'Fetch.AI"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 3)
$c.Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.Value = "'2.22"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 5)
$c.Value = "'  -9.10%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 4)
$c.Value = "'10.04"
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 5)
$c.Value = "'  -6.29%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 5)
$c.Value = "'  -0.05%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 5)
$c.Value = "'  -0.05%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 5)
$c.Value = "'  -2.54%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 5)
$c.Value = "'  -4.55%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 4)
$c.Value = "'2.09"
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 5)
$c.Value = "'  -5.97%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 4)
$c.Value = "'28.08"
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 5)
$c.Value = "'  +2.39%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.109"
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 5)
$c.Value = "'  -2.14%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.0₃0810"
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 5)
$c.Value = "'  -1.75%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 5)
$c.Value = "'  -3.13%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 4)
$c.Value = "'5.77"
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 5)
$c.Value = "'  -3.68%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 4)
$c.Value = "'2.12"
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 5)
$c.Value = "'  -5.37%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 4)
$c.Value = "'50.29"
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 5)
$c.Value = "'  -0.68%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.Value = "'9.15"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 5)
$c.Value = "'  -1.45%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.Value = "'2.90"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 5)
$c.Value = "'  -12.98%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 5)
$c.Value = "'  +4.14%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.Value = "'392.44"
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 5)
$c.Value = "'  -10.83%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.Value = "'0.0359"
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 5)
$c.Value = "'  -1.68%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.Value = "'0.269"
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 5)
$c.Value = "'  -8.04%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 4)
$c.Value = "'2.736.45"
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 5)
$c.Value = "'  -3.95%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 4)
$c.Value = "'37.27"
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 5)
$c.Value = "'  -4.23%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.Value = "'129.44"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 5)
$c.Value = "'  -0.74%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 5)
$c.Value = "'  +0.05%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 5)
$c.Value = "'  -1.49%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 5)
$c.Value = "'  -0.96%  "
$c.Style = "Normal"
